$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22: the phone number in A22 was stored as text; normalize it to a
#     plain number (matches the rest of the "phone" column). ---
$ws.Range("A22").Value = 71277628

# --- New row 23: payment 71277628 (Cash) 2025-08-18T16:53:40 ---
# A23 keeps the phone number as text (leading apostrophe forces text entry
# without leaving the cell visually/format-wise different from a plain
# General cell).
$ws.Range("A23").Formula = "'71277628"
$ws.Range("A23").Style = "Normal"

$ws.Range("B23").Formula = "'"
$ws.Range("B23").Style = "Normal"

$ws.Range("C23").Value = "Cash"
$ws.Range("D23").Value = "2025-08-18T16:53:40"
$ws.Range("E23").Value = 760

$ws.Range("F23").Formula = "'"
$ws.Range("F23").Style = "Normal"

$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 760
